$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1711229946524064
$ws.Range("C2").Value = 0.6283422459893048
$ws.Range("J2").Value = 0.0160427807486631
$ws.Range("P2").Value = 0.1310160427807487
$ws.Range("S2").Value = 0.053475935828877
$ws.Range("B3").Value = 0.007843137254901961
$ws.Range("C3").Value = 0.05882352941176471
$ws.Range("J3").Value = 0.02352941176470588
$ws.Range("O3").Value = 0.00392156862745098
$ws.Range("P3").Value = 0.7607843137254902
$ws.Range("S3").Value = 0.1450980392156863
$ws.Range("J4").Value = 0.03636363636363636
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.3272727272727273
$ws.Range("B6").Value = 0.04017857142857143
$ws.Range("D6").Value = 0.01339285714285714
$ws.Range("F6").Value = 0.04910714285714286
$ws.Range("J6").Value = 0.2901785714285715
$ws.Range("O6").Value = 0.04017857142857143
$ws.Range("Q6").Value = 0.1339285714285714
$ws.Range("R6").Value = 0.1026785714285714
$ws.Range("S6").Value = 0.3303571428571428
$ws.Range("B7").Value = 0.1538461538461539
$ws.Range("D7").Value = 0.02403846153846154
$ws.Range("F7").Value = 0.03846153846153846
$ws.Range("J7").Value = 0.1971153846153846
$ws.Range("O7").Value = 0.03365384615384615
$ws.Range("Q7").Value = 0.1490384615384615
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.3269230769230769
$ws.Range("B8").Value = 0.1180400890868597
$ws.Range("D8").Value = 0.0311804008908686
$ws.Range("E8").Value = 0.0022271714922049
$ws.Range("F8").Value = 0.06013363028953229
$ws.Range("J8").Value = 0.1358574610244989
$ws.Range("O8").Value = 0.0178173719376392
$ws.Range("Q8").Value = 0.2160356347438753
$ws.Range("R8").Value = 0.1158129175946548
$ws.Range("S8").Value = 0.3028953229398664
$ws.Range("B9").Value = 0.1372549019607843
$ws.Range("D9").Value = 0.0261437908496732
$ws.Range("F9").Value = 0.07843137254901961
$ws.Range("J9").Value = 0.08496732026143791
$ws.Range("O9").Value = 0.03267973856209151
$ws.Range("Q9").Value = 0.2418300653594771
$ws.Range("R9").Value = 0.08496732026143791
$ws.Range("S9").Value = 0.3137254901960784
$ws.Range("B10").Value = 0.1471048513302035
$ws.Range("D10").Value = 0.02269170579029734
$ws.Range("E10").Value = 0.001564945226917058
$ws.Range("F10").Value = 0.07511737089201878
$ws.Range("J10").Value = 0.1517996870109546
$ws.Range("O10").Value = 0.01721439749608764
$ws.Range("Q10").Value = 0.1964006259780908
$ws.Range("R10").Value = 0.08528951486697965
$ws.Range("S10").Value = 0.3028169014084507
$ws.Range("G11").Value = 0.1601307189542484
$ws.Range("J11").Value = 0.09477124183006536
$ws.Range("K11").Value = 0.2189542483660131
$ws.Range("L11").Value = 0.5196078431372549
$ws.Range("S11").Value = 0.006535947712418301
$ws.Range("G12").Value = 0.7901234567901234
$ws.Range("J12").Value = 0.154320987654321
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.01851851851851852
$ws.Range("S12").Value = 0.02469135802469136
$ws.Range("G13").Value = 0.7454545454545455
$ws.Range("J13").Value = 0.2363636363636364
$ws.Range("S13").Value = 0.01818181818181818
$ws.Range("F15").Value = 0.01015228426395939
$ws.Range("H15").Value = 0.16751269035533
$ws.Range("I15").Value = 0.04060913705583756
$ws.Range("J15").Value = 0.3350253807106599
$ws.Range("K15").Value = 0.04060913705583756
$ws.Range("M15").Value = 0.02538071065989848
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.3248730964467005
$ws.Range("F16").Value = 0.04089219330855019
$ws.Range("H16").Value = 0.1970260223048327
$ws.Range("I16").Value = 0.06691449814126393
$ws.Range("J16").Value = 0.4163568773234201
$ws.Range("K16").Value = 0.09293680297397769
$ws.Range("M16").Value = 0.02973977695167286
$ws.Range("N16").Value = 0.003717472118959108
$ws.Range("O16").Value = 0.03345724907063197
$ws.Range("S16").Value = 0.1189591078066914
$ws.Range("F17").Value = 0.02031602708803612
$ws.Range("H17").Value = 0.1941309255079007
$ws.Range("I17").Value = 0.08126410835214447
$ws.Range("J17").Value = 0.4243792325056434
$ws.Range("K17").Value = 0.09932279909706546
$ws.Range("M17").Value = 0.02031602708803612
$ws.Range("O17").Value = 0.05417607223476298
$ws.Range("S17").Value = 0.1060948081264108
$ws.Range("F18").Value = 0.01941747572815534
$ws.Range("H18").Value = 0.2087378640776699
$ws.Range("I18").Value = 0.07281553398058252
$ws.Range("J18").Value = 0.383495145631068
$ws.Range("K18").Value = 0.1116504854368932
$ws.Range("M18").Value = 0.02912621359223301
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.04854368932038835
$ws.Range("S18").Value = 0.1213592233009709
$ws.Range("F19").Value = 0.01717557251908397
$ws.Range("H19").Value = 0.2194656488549618
$ws.Range("I19").Value = 0.07538167938931298
$ws.Range("J19").Value = 0.3740458015267176
$ws.Range("K19").Value = 0.1297709923664122
$ws.Range("M19").Value = 0.02862595419847328
$ws.Range("N19").Value = 0.001908396946564885
$ws.Range("O19").Value = 0.06202290076335878
$ws.Range("S19").Value = 0.0916030534351145
